$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.217309109643015
$ws.Range("C2").Value = 4.895899660347447
$ws.Range("D2").Value = 5.219366465294583
$ws.Range("F2").Value = 27.06372233065336
$ws.Range("G2").Value = 33.76926825462541
$ws.Range("H2").Value = 15.35851858051627
$ws.Range("K2").Value = 8.53723674707112
$ws.Range("M2").Value = 19.97530439811406
$ws.Range("N2").Value = 19.3691088453517

$ws.Range("B3").Value = 8.981433454413791
$ws.Range("C3").Value = 4.774958898174502
$ws.Range("D3").Value = 5.196141735060422
$ws.Range("F3").Value = 26.99952462110661
$ws.Range("G3").Value = 33.62950023259593
$ws.Range("H3").Value = 15.38434201346014
$ws.Range("K3").Value = 8.389330491237031
$ws.Range("M3").Value = 19.38202190400647
$ws.Range("N3").Value = 19.42977637976933

$ws.Range("B4").Value = 8.835786290094894
$ws.Range("C4").Value = 4.698123739034182
$ws.Range("D4").Value = 5.181455967501355
$ws.Range("F4").Value = 26.96703258262055
$ws.Range("G4").Value = 33.55352056854566
$ws.Range("H4").Value = 15.40282289812294
$ws.Range("K4").Value = 8.299225291041834
$ws.Range("M4").Value = 19.01659802019375
$ws.Range("N4").Value = 19.46888453440714

$ws.Range("B5").Value = 8.776327097726636
$ws.Range("C5").Value = 4.666188623533321
$ws.Range("D5").Value = 5.175365702882949
$ws.Range("F5").Value = 26.95554054449997
$ws.Range("G5").Value = 33.52505171708774
$ws.Range("H5").Value = 15.41101318196059
$ws.Range("K5").Value = 8.262737204131199
$ws.Range("M5").Value = 18.86764624249121
$ws.Range("N5").Value = 19.48528959959164

$ws.Range("B6").Value = 8.76645046888712
$ws.Range("C6").Value = 4.660848855290743
$ws.Range("D6").Value = 5.174348074997267
$ws.Range("F6").Value = 26.95373809411303
$ws.Range("G6").Value = 33.520475582966
$ws.Range("H6").Value = 15.41241295496094
$ws.Range("K6").Value = 8.256693832079854
$ws.Range("M6").Value = 18.84291810217589
$ws.Range("N6").Value = 19.48804195294761

$ws.Range("B7").Value = 8.834984701109599
$ws.Range("C7").Value = 4.697695545332333
$ws.Range("D7").Value = 5.18137425825585
$ws.Range("F7").Value = 26.96687050836189
$ws.Range("G7").Value = 33.55312650880376
$ws.Range("H7").Value = 15.40293068757269
$ws.Range("K7").Value = 8.298732197689279
$ws.Range("M7").Value = 19.01458902241837
$ws.Range("N7").Value = 19.46910388203826

$ws.Range("B8").Value = 9.136207166307548
$ws.Range("C8").Value = 4.854747916808733
$ws.Range("D8").Value = 5.211447731707866
$ws.Range("F8").Value = 27.04015443768767
$ws.Range("G8").Value = 33.71904709552842
$ws.Range("H8").Value = 15.36687709599117
$ws.Range("K8").Value = 8.486122448185613
$ws.Range("M8").Value = 19.77112443080883
$ws.Range("N8").Value = 19.38964195081922

$ws.Range("B9").Value = 9.716287808207325
$ws.Range("C9").Value = 5.141321731131037
$ws.Range("D9").Value = 5.266986654790154
$ws.Range("F9").Value = 27.23841616172269
$ws.Range("G9").Value = 34.12145164362891
$ws.Range("H9").Value = 15.31704563723854
$ws.Range("K9").Value = 8.856965611006476
$ws.Range("M9").Value = 21.23529182406624
$ws.Range("N9").Value = 19.24851564455164

$ws.Range("B10").Value = 10.13085177429639
$ws.Range("C10").Value = 5.337659789279462
$ws.Range("D10").Value = 5.305622603666517
$ws.Range("F10").Value = 27.41668317894058
$ws.Range("G10").Value = 34.46244362415755
$ws.Range("H10").Value = 15.29320730794206
$ws.Range("K10").Value = 9.128518200917576
$ws.Range("M10").Value = 22.28646077334837
$ws.Range("N10").Value = 19.15372719942405

$ws.Range("B11").Value = 10.31595003317105
$ws.Range("C11").Value = 5.42368256956825
$ws.Range("D11").Value = 5.322712771262961
$ws.Range("F11").Value = 27.50468454786051
$ws.Range("G11").Value = 34.62700113022541
$ws.Range("H11").Value = 15.28514471673489
$ws.Range("K11").Value = 9.25125902457056
$ws.Range("M11").Value = 22.75696930947696
$ws.Range("N11").Value = 19.11252402755295

$ws.Range("B12").Value = 10.3854673134243
$ws.Range("C12").Value = 5.455768114098496
$ws.Range("D12").Value = 5.329113231945664
$ws.Range("F12").Value = 27.53898385489197
$ws.Range("G12").Value = 34.690632313763
$ws.Range("H12").Value = 15.28249208239367
$ws.Range("K12").Value = 9.297576857429572
$ws.Range("M12").Value = 22.93386134027847
$ws.Range("N12").Value = 19.09719611726489

$ws.Range("B13").Value = 10.37052222326885
$ws.Range("C13").Value = 5.448879923190409
$ws.Range("D13").Value = 5.327737967585082
$ws.Range("F13").Value = 27.53155382325802
$ws.Range("G13").Value = 34.67687030669521
$ws.Range("H13").Value = 15.2830455558904
$ws.Range("K13").Value = 9.287609402557282
$ws.Range("M13").Value = 22.89582399503377
$ws.Range("N13").Value = 19.10048504549824

$ws.Range("B14").Value = 10.32168119414064
$ws.Range("C14").Value = 5.426332161371193
$ws.Range("D14").Value = 5.323240776214702
$ws.Range("F14").Value = 27.50748693957509
$ws.Range("G14").Value = 34.63220995623404
$ws.Range("H14").Value = 15.28491845398027
$ws.Range("K14").Value = 9.255073122608035
$ws.Range("M14").Value = 22.77154891009991
$ws.Range("N14").Value = 19.11125748732472

$ws.Range("B15").Value = 10.29168763515643
$ws.Range("C15").Value = 5.41245682161103
$ws.Range("D15").Value = 5.320476801171075
$ws.Range("F15").Value = 27.49287169977303
$ws.Range("G15").Value = 34.6050244444672
$ws.Range("H15").Value = 15.28611782707783
$ws.Range("K15").Value = 9.235121302786336
$ws.Range("M15").Value = 22.69525529522716
$ws.Range("N15").Value = 19.1178916884835

$ws.Range("B16").Value = 10.11867852002118
$ws.Range("C16").Value = 5.331970346739405
$ws.Range("D16").Value = 5.304495831253697
$ws.Range("F16").Value = 27.41106955671495
$ws.Range("G16").Value = 34.45187594079827
$ws.Range("H16").Value = 15.29379023596244
$ws.Range("K16").Value = 9.120476763203031
$ws.Range("M16").Value = 22.25554209649647
$ws.Range("N16").Value = 19.15645842415999

$ws.Range("B17").Value = 10.01159502956225
$ws.Range("C17").Value = 5.281739390323835
$ws.Range("D17").Value = 5.29456667077324
$ws.Range("F17").Value = 27.36264325141697
$ws.Range("G17").Value = 34.36031387769743
$ws.Range("H17").Value = 15.29920984094078
$ws.Range("K17").Value = 9.049909002752129
$ws.Range("M17").Value = 21.98369633148512
$ws.Range("N17").Value = 19.18060819710647

$ws.Range("B18").Value = 9.949679054925717
$ws.Range("C18").Value = 5.252538833168711
$ws.Range("D18").Value = 5.288810063435566
$ws.Range("F18").Value = 27.33544101477883
$ws.Range("G18").Value = 34.30854080931621
$ws.Range("H18").Value = 15.3025888481484
$ws.Range("K18").Value = 9.009248683454215
$ws.Range("M18").Value = 21.82662499388166
$ws.Range("N18").Value = 19.19467894870405

$ws.Range("B19").Value = 9.928662003195114
$ws.Range("C19").Value = 5.242599460168465
$ws.Range("D19").Value = 5.286853183227244
$ws.Range("F19").Value = 27.3263431888434
$ws.Range("G19").Value = 34.2911655544914
$ws.Range("H19").Value = 15.30377786772148
$ws.Range("K19").Value = 8.995470993722504
$ws.Range("M19").Value = 21.77332654662257
$ws.Range("N19").Value = 19.19947407784141

$ws.Range("B20").Value = 10.02302837577395
$ws.Range("C20").Value = 5.28711866634044
$ws.Range("D20").Value = 5.295628378205868
$ws.Range("F20").Value = 27.36773103367936
$ws.Range("G20").Value = 34.36996888351251
$ws.Range("H20").Value = 15.29860581601784
$ws.Range("K20").Value = 9.057428826540031
$ws.Range("M20").Value = 22.01270988237681
$ws.Range("N20").Value = 19.17801874319435

$ws.Range("B21").Value = 10.33604316181844
$ws.Range("C21").Value = 5.432968387992535
$ws.Range("D21").Value = 5.324563654137835
$ws.Range("F21").Value = 27.51452965518185
$ws.Range("G21").Value = 34.6452923813528
$ws.Range("H21").Value = 15.28435746572616
$ws.Range("K21").Value = 9.264634568989885
$ws.Range("M21").Value = 22.80808755323138
$ws.Range("N21").Value = 19.10808590703155

$ws.Range("B22").Value = 10.53723259284068
$ws.Range("C22").Value = 5.525430564285868
$ws.Range("D22").Value = 5.343058798633841
$ws.Range("F22").Value = 27.61614512125516
$ws.Range("G22").Value = 34.83288702624247
$ws.Range("H22").Value = 15.27737987356204
$ws.Range("K22").Value = 9.399093832042205
$ws.Range("M22").Value = 23.32038622368546
$ws.Range("N22").Value = 19.06398237314809

$ws.Range("B23").Value = 10.43018639554793
$ws.Range("C23").Value = 5.476348254392343
$ws.Range("D23").Value = 5.33322608866601
$ws.Range("F23").Value = 27.56139825536871
$ws.Range("G23").Value = 34.73207787626031
$ws.Range("H23").Value = 15.28089020260372
$ws.Range("K23").Value = 9.327433617235185
$ws.Range("M23").Value = 23.04770519167154
$ws.Range("N23").Value = 19.08737495850187

$ws.Range("B24").Value = 10.0178604514716
$ws.Range("C24").Value = 5.284687697729735
$ws.Range("D24").Value = 5.295148530304101
$ws.Range("F24").Value = 27.36542885658566
$ws.Range("G24").Value = 34.36560115127495
$ws.Range("H24").Value = 15.29887807590326
$ws.Range("K24").Value = 9.05402939079716
$ws.Range("M24").Value = 21.99959528318082
$ws.Range("N24").Value = 19.17918885364213

$ws.Range("B25").Value = 9.561050761755732
$ws.Range("C25").Value = 5.066213712236584
$ws.Range("D25").Value = 5.252337385847698
$ws.Range("F25").Value = 27.17900099224918
$ws.Range("G25").Value = 34.0044966900619
$ws.Range("H25").Value = 15.32828650364864
$ws.Range("K25").Value = 8.756592622964893
$ws.Range("M25").Value = 20.8426488938148
$ws.Range("N25").Value = 19.2851271173574

Write-Output "Done updating loading_percent values"